# The underlying data rows (2-12) were reshuffled: each row now shows the
# Fecha/Calidad/Volumen/Precio* values that used to belong to a different row.
# Only the cells whose value actually differs from before are written below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fecha/Calidad/Volumen/Precio* now come from original row 6
$ws.Range("D2").Value = 44473
$ws.Range("M2").Value = 180

# Row 3: Fecha/Calidad/Volumen/Precio* now come from original row 2
$ws.Range("D3").Value = 44435
$ws.Range("M3").Value = 40

# Row 4: Fecha/Calidad/Volumen/Precio* now come from original row 9
$ws.Range("D4").Value = 44517
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 27000
$ws.Range("O4").Value = 27000
$ws.Range("P4").Value = 27000
$ws.Range("S4").Value = 2700

# Row 5: Fecha/Calidad/Volumen/Precio* now come from original row 10
$ws.Range("D5").Value = 44517
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 2500

# Row 6: Fecha/Calidad/Volumen/Precio* now come from original row 3
$ws.Range("D6").Value = 44434
$ws.Range("M6").Value = 20

# Row 7: Fecha/Calidad/Volumen/Precio* now come from original row 11
$ws.Range("D7").Value = 44432
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 2000

# Row 8: Fecha/Calidad/Volumen/Precio* now come from original row 7
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 30000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 30000
$ws.Range("S8").Value = 3000

# Row 9: Fecha/Calidad/Volumen/Precio* now come from original row 8
$ws.Range("D9").Value = 44503
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 25000
$ws.Range("S9").Value = 2500

# Row 10: Fecha/Calidad/Volumen/Precio* now come from original row 12
$ws.Range("D10").Value = 44476
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 2000

# Row 11: Fecha/Calidad/Volumen/Precio* now come from original row 4
$ws.Range("D11").Value = 44466
$ws.Range("M11").Value = 60

# Row 12: Fecha/Calidad/Volumen/Precio* now come from original row 5
$ws.Range("D12").Value = 44511
$ws.Range("N12").Value = 28000
$ws.Range("O12").Value = 28000
$ws.Range("P12").Value = 28000
$ws.Range("S12").Value = 2800
